$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new "line" entries (line7, line8) were inserted into the name pool
# right after line6. Because of that, rows 8-15 (which kept referencing the
# same relative slot in the name list) now show the next names down the
# list, and two brand new rows (16, 17) are appended holding the final two
# "extr" names. Re-point/refresh every name cell plus the changed C/D/E
# values so the sheet ends up matching the target table below:
#
# row  A   B(name)  C   D   E
#  2   0   line1    7   9  TRUE
#  3   1   line2    9   8  TRUE
#  4   2   line3    8  10  FALSE
#  5   3   line4    8  11  TRUE
#  6   4   line5   10   5  TRUE
#  7   5   line6   12   8  TRUE
#  8   6   line7   14  11  TRUE
#  9   7   line8   16   9  TRUE
# 10   8   extr1    5  12  TRUE
# 11   9   extr2    5   9  TRUE
# 12  10   extr3   10  11  FALSE
# 13  11   extr4    7   8  FALSE
# 14  12   extr5    9  11  TRUE
# 15  13   extr6    7  11  TRUE
# 16  14   extr7    5   7  FALSE
# 17  15   extr8    8   5  FALSE

# --- Rows 8 and 9: name shifts to line7 / line8, plus C/D/E updates ---
$ws.Range("B8").Value = "line7"
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $true

$ws.Range("B9").Value = "line8"
$ws.Range("C9").Value = 16
$ws.Range("D9").Value = 9
$ws.Range("E9").Value = $true

# --- Rows 10-15: name shifts down to the next "extr" entry ---
$ws.Range("B10").Value = "extr1"
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12
$ws.Range("E10").Value = $true

$ws.Range("B11").Value = "extr2"
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9
$ws.Range("E11").Value = $true

$ws.Range("B12").Value = "extr3"
$ws.Range("C12").Value = 10
$ws.Range("D12").Value = 11
$ws.Range("E12").Value = $false

$ws.Range("B13").Value = "extr4"
$ws.Range("C13").Value = 7
$ws.Range("D13").Value = 8
$ws.Range("E13").Value = $false

$ws.Range("B14").Value = "extr5"
$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11
$ws.Range("E14").Value = $true

$ws.Range("B15").Value = "extr6"
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11
$ws.Range("E15").Value = $true

# --- New rows 16 and 17, formatted like the existing data rows ---
$ws.Range("A15").Copy()
$ws.Range("A16:A17").PasteSpecial(-4122)
$ws.Range("B15:E15").Copy()
$ws.Range("B16:E17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "extr7"
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = $false

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "extr8"
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = $false
